# Append a new log row (row 35) to Sheet1, mirroring the existing
# "SKIPPED" run-log entries. New row reuses the same cell styling as the
# previous last row (row 34) via a copy/paste-format, then the actual
# values for the new run are written in.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy formatting (style) of the last existing data row onto the new row.
$ws.Range("A34:H34").Copy()
$ws.Range("A35:H35").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# Fill in the new run-log entry.
$ws.Cells.Item(35, 1).Value = "2025-08-20 06:49:59 UTC"
$ws.Cells.Item(35, 2).Value = "2025-08-20 12:19:59 IST"
$ws.Cells.Item(35, 3).Value = "SKIPPED"
$ws.Cells.Item(35, 4).Value = "No change in PDF. Skipping download & Excel update."
$ws.Cells.Item(35, 5).Value = "https://nalcoindia.com/wp-content/uploads/2025/08/INGOT-15-08-2025.pdf"
$ws.Cells.Item(35, 6).Value = ""
$ws.Cells.Item(35, 7).Value = 0
$ws.Cells.Item(35, 8).Value = ""
